# Dispatch.xlsx edit: add Items sheet with item id/description data,
# pulled from a "request Item data from file" workflow and uploaded into a
# new worksheet; update Dispatch header H1 to "itemDescription" and refresh
# selections/active tab to reflect the new sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Dispatch sheet: rename column header H1 from "description" to
#    "itemDescription", and move the remembered selection off the sheet.
# ---------------------------------------------------------------------------
$dispatch = $wb.Worksheets.Item("Dispatch")
$dispatch.Range("H1").Value = "itemDescription"
$dispatch.Activate()
$dispatch.Range("A28").Select()

# ---------------------------------------------------------------------------
# 2) Add new "Items" worksheet after "Resources", populated with the item
#    catalog (_id / description) pulled from the source file.
# ---------------------------------------------------------------------------
    @("03723", "PREMIUM GD BTM RUBBER 9' BLACK RL P0960"),
    @("03749", "PREMIUM GD BTM RUBBER 16' BLACK RL P0960"),
    @("03822", "GARAGE DOOR SEAL F/TOP&SIDES 30' WHIT..."),
    @("05991", "U-SHAPED DR BTM VINYL 36`" BROWN LV404..."),
    @("06528", "U-SHAPED DR BTM VINYL 36`" WHITE LV404..."),
    @("08460", "REPLACE VINYL UNIV GD BTM 9' BLACK P0199"),
    @("08462", "REPLACE VINYL UNIV GD BTM 18' BLACK P..."),
    @("09060", "PLATINUM REPLACE VINYL UNIV GD BTM 18..."),
    @("23647", "WALLBASE 4`"X48`" ADHES BEIGE V3951"),
    @("29520059", "GARAGE DB EPDM RUBBER 100' RL BLACK 2..."),
    @("29520098", "UNIVERSAL GD BTM 9' D#1729 P0199 2RRT5"),
    @("29520127", "WALLBASE 4`"X60' ALMOND BULK V3951 2R..."),
    @("29520128", "WALLBASE 4`"X60' BLACK BULK V3951 2RRW9"),
    @("29520129", "WALLBASE 4`"X60' BROWN BULK V3951 2RRX1"),
    @("29520131", "WALLBASE 4`"X60' WHITE BULK V3951 2RRX3"),
    @("29520136", "VINYL STAIR NOSING 36`" BLACK #1716 P0..."),
    @("29520215", "GARAGE DOOR THRESHOLD KIT 20' GRAY 5M..."),
    @("29520230", "VINYL WALL BASE 4`"X4' BLACK 5MFJ6 V3951"),
    @("29520233", "VINYL WALL BASE 4`"X4' SNOW WHITE 5MFJ..."),
    @("29702", "VINYL STAIR NOSING 36`" BLACK #1716 P0919"),
    @("29710", "VINYL STAIR NOSING 36`" BROWN #1716 P0919"),
    @("31147", "CA UNIVERSAL GD BTM 16' BLACK D#1729 ..."),
    @("31391", "DIVIDER T #11 SHORT EDGE 144`" BLACK P..."),
    @("38339", "GARAGE DOOR SEAL TOP & SIDES 30' WHIT..."),
    @("38342", "BOTTOM GARAGE DOOR SEAL BLACK 16' P0960"),
    @("38386", "PLATINUM COLLECTION DOOR WEATHERSTRIP..."),
    @("38387", "PLATINUM COLLECTION DOOR WEATHERSTRIP..."),
    @("38388", "PLATINUM COLLECTION DOOR WEATHERSTRIP..."),
    @("40200", "ST CINCH SLIDE ON UNDER DR SEAL CLR 3..."),
    @("40717", "DENY SLIDE ON UNDER DOOR SEAL WHITE 3..."),
    @("40718", "DENY SLIDE ON UNDER DOOR SEAL BROWN 3..."),
    @("43336", "CINCH SLIDE ON UNDER DOOR SEAL WHITE ..."),
    @("43337", "CINCH SLIDE ON UNDER DOOR SEAL BROWN ..."),
    @("43339", "CINCH SLIDE ON UNDER DOOR SEAL BLACK ..."),
    @("43346", "PLATINUM COLLECTION DOOR W/S PREM REP..."),
    @("43347", "PLATINUM COLLECTION DOOR W/S PREM REP..."),
    @("43348", "PLATINUM COLLECTION DOOR W/S PREM REP..."),
    @("43398", "LIGHT & SOUND UNDER DOOR SEAL 32`" WHITE"),
    @("43399", "LIGHT & SOUND UNDER DOOR SEAL 32`" BROWN"),
    @("50101", "THRESHOLD GARAGE DOOR VINYL 20' P0955"),
    @("70056", "VINYL CLAD 96`" BROWN LV9650"),
    @("70284040", "WIP-VINYL CLAD BRONZE 84`" NO/LBL LV9650"),
    @("70284047", "WIP-VINYL CLAD WHITE 36`" NO/LBL LV9650"),
    @("73896", "WALLBASE 4`"X60' BLACK BULK V3951"),
    @("73897", "WALLBASE 4`"X60' WHITE BULK V3951"),
    @("73898", "WALLBASE 4`"X60' SILVER GRAY BULK V3951"),
    @("73899", "WALLBASE 4`"X60' ALMOND BULK V3951"),
    @("73900", "WALLBASE 4`"X60' BROWN BULK V3951"),
    @("75317", "WALLBASE 4`"X48`" SNOW WHITE V3951"),
    @("75457", "WALLBASE 4`"X120' BLACK BULK V3951"),
    @("75465", "WALLBASE 4`"X120' BROWN BULK V3951"),
    @("75507", "WALLBASE 4`"X120' WHITE BULK V3951"),
    @("75556", "RESIDENT STAIR TREAD 24`" BLACK P0047"),
    @("75572", "RESIDENT STAIR TREAD 24`" BROWN P0047"),
    @("75598", "WALLBASE 2-1/2X48`" BLACK P1032"),
    @("75697", "WALLBASE 2-1/2X48`" WHITE P1032"),
    @("75903", "WALLBASE 2-1/2`"X120' BROWN BULK P1032"),
    @("78899", "ULTRA DR BTM DB347 36`" BROWN LV5219 V..."),
    @("82586", "DR BTM DB346 36`" BROWN LV3871 V1346"),
    @("84378196", "VINYL DR GASKET 12`" BRN P0226"),
    @("84378394", "VINYL DR GASKET 12`" WHT P0226"),
    @("84382831", "TPV `"D`" PROFILE 12`" BROWN V4415"),
    @("87643", "UNIVERSAL GD BTM 9' BLACK D#1729 P0199"),
    @("87668", "UNIVERSAL GD BTM 18' BLACK D#1729 P0199"),
    @("88660001", "VINYL DIVIDER T #9 BLACK 97`" P0926"),
    @("88660002", "VINYL DIVIDER T #9 LT BROWN 97`" P0926"),
    @("91396", "VINYL CLAD REPLACE 96`" WHITE LV9650"),
    @("91397", "VINYL CLAD REPLACE 96`" BROWN LV9650"),
    @("91800", "VINYL CLAD REPLACE 81`" BEIGE LV9650"),
    @("91868", "VINYL CLAD REPLACE 81`" WHITE LV9650"),
    @("93229", "WALLBASE 4`"X20' ADHES ALMOND V3951"),
    @("CF20111", "REPL. FOR GARAGE DOOR BOT. 16' DENY B..."),
    @("CF20644", "DENY CINCH SLIDE ON UNDER DOOR SEAL W..."),
    @("CF20645", "DENY CINCH SLIDE ON UNDER DOOR SEAL B..."),
    @("CF20920", "PREMIUM RUBBER NAIL-ON GARAGE DOOR BO..."),
    @("CF21079", "GARAGE DOOR BOT. 1-5/8'' x 9' ALU + V..."),
    @("CF21180", "UNIVERSAL GARAGE DOOR BTM 16' BLACK D..."),
    @("CF28416", "CLIMALOC PLUS COLLECTION DOOR W/S PRE..."),
    @("CF65821", "BULK REPL. KERF 81'' PERMA-FOAM WHITE..."),
    @("CF65822", "BULK REPL. KERF 81'' PERMA-FOAM BLACK..."),
    @("CFRP278", "REPL. FOR GARAGE DOOR BOT. 9FT BLACK"),
    @("VD004153BZ084", "PREMIUM DOOR SEAL 84`" BROWN"),
    @("VD004153WH084", "PREMIUM DOOR SEAL 84`" WHITE"),
    @("WIP004153BZ036", "PLATINUM BRONZE 36`" LV4153 NO LABEL"),
    @("WIP004153BZ084", "PLATINUM BRONZE 84`" LV4153 NO LABEL"),
    @("WIP004153WH036", "PLATINUM WHITE 36`" LV4153 NO LABEL"),
    @("WIP004153WH084", "PLATINUM WHITE 84`" LV4153 NO LABEL"),
    @("WIP27950239", "VINYL INSERT SDS 96`" GRAY P0013")

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$items = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$items.Name = "Items"

# Header row: reuse the existing bold-with-bottom-border look from the
# Dispatch header, recolored to match the Resources header's font.
$items.Range("A1").Value = "_id"
$items.Range("B1").Value = "description"
$dispatch.Range("A1").Copy()
$items.Range("A1:B1").PasteSpecial(-4122)
$items.Range("A1:B1").Font.Color = 0

# Data rows. Column A holds item numbers that must stay text (leading
# zeros, etc.) -- format the column as Text before writing, then drop back
# to the workbook's default style so the saved cells carry no explicit
# number format (matches how the source data was typed in).
$dataRange = $items.Range("A2:B" + (1 + $itemsData.Count))
$items.Range("A2:A" + (1 + $itemsData.Count)).NumberFormat = "@"

for ($i = 0; $i -lt $itemsData.Count; $i++) {
    $row = $i + 2
    $items.Cells.Item($row, 1).Value = $itemsData[$i][0]
    $items.Cells.Item($row, 2).Value = $itemsData[$i][1]
}

$dataRange.Style = "Normal"

$items.Columns.Item(2).AutoFit()

$items.Activate()
$items.Range("B16").Select()
